# Add team record (Wins/Losses/Ties) columns to the player data sheet.
# The source data previously spanned A1:AC53; we extend it with three
# new columns (AD, AE, AF) holding the team's win/loss/tie record,
# repeated for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching the style of the existing header row (A1:AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Team record is constant for every player row (rows 2-53).
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 73  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 89  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
